$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gip"
$ws.Cells.Item(2, 3).Value = "Dpp4"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.01853566666666667
$ws.Cells.Item(2, 8).Value = 0.055607
$ws.Cells.Item(2, 9).Value = 0.04056456431044909
$ws.Cells.Item(2, 10).Value = 0.04056456431044909
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 4.500364333333334
$ws.Cells.Item(2, 14).Value = 13.501093
$ws.Cells.Item(2, 15).Value = 0.7110918985538353
$ws.Cells.Item(2, 16).Value = 0.7110918985538353
$ws.Cells.Item(2, 17).Value = 0.08341725316122223
$ws.Cells.Item(2, 18).Value = 0.750755278451
$ws.Cells.Item(2, 19).Value = 0.02884513304952639
$ws.Cells.Item(2, 20).Value = 0.02884513304952639

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gip"
$ws.Cells.Item(3, 3).Value = "Dpp4"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.01853566666666667
$ws.Cells.Item(3, 8).Value = 0.055607
$ws.Cells.Item(3, 9).Value = 0.04056456431044909
$ws.Cells.Item(3, 10).Value = 0.04056456431044909
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.116902
$ws.Cells.Item(3, 14).Value = 0.350706
$ws.Cells.Item(3, 15).Value = 0.0184714078611429
$ws.Cells.Item(3, 16).Value = 0.01847140786114289
$ws.Cells.Item(3, 17).Value = 0.002166856504666667
$ws.Cells.Item(3, 18).Value = 0.019501708542
$ws.Cells.Item(3, 19).Value = 0.000749284612087866
$ws.Cells.Item(3, 20).Value = 0.0007492846120878658

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Gip"
$ws.Cells.Item(4, 3).Value = "Dpp4"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.01853566666666667
$ws.Cells.Item(4, 8).Value = 0.055607
$ws.Cells.Item(4, 9).Value = 0.04056456431044909
$ws.Cells.Item(4, 10).Value = 0.04056456431044909
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.711542
$ws.Cells.Item(4, 14).Value = 5.134626
$ws.Cells.Item(4, 15).Value = 0.2704366935850219
$ws.Cells.Item(4, 16).Value = 0.2704366935850219
$ws.Cells.Item(4, 17).Value = 0.031724571998
$ws.Cells.Item(4, 18).Value = 0.285521147982
$ws.Cells.Item(4, 19).Value = 0.01097014664883484
$ws.Cells.Item(4, 20).Value = 0.01097014664883484

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gip"
$ws.Cells.Item(5, 3).Value = "Dpp4"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.4384066666666667
$ws.Cells.Item(5, 8).Value = 1.31522
$ws.Cells.Item(5, 9).Value = 0.9594354356895509
$ws.Cells.Item(5, 10).Value = 0.9594354356895509
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 4.500364333333334
$ws.Cells.Item(5, 14).Value = 13.501093
$ws.Cells.Item(5, 15).Value = 0.7110918985538353
$ws.Cells.Item(5, 16).Value = 0.7110918985538353
$ws.Cells.Item(5, 17).Value = 1.972989726162222
$ws.Cells.Item(5, 18).Value = 17.75690753546
$ws.Cells.Item(5, 19).Value = 0.6822467655043088
$ws.Cells.Item(5, 20).Value = 0.6822467655043088

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gip"
$ws.Cells.Item(6, 3).Value = "Dpp4"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.4384066666666667
$ws.Cells.Item(6, 8).Value = 1.31522
$ws.Cells.Item(6, 9).Value = 0.9594354356895509
$ws.Cells.Item(6, 10).Value = 0.9594354356895509
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.116902
$ws.Cells.Item(6, 14).Value = 0.350706
$ws.Cells.Item(6, 15).Value = 0.0184714078611429
$ws.Cells.Item(6, 16).Value = 0.01847140786114289
$ws.Cells.Item(6, 17).Value = 0.05125061614666666
$ws.Cells.Item(6, 18).Value = 0.46125554532
$ws.Cells.Item(6, 19).Value = 0.01772212324905503
$ws.Cells.Item(6, 20).Value = 0.01772212324905503

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Gip"
$ws.Cells.Item(7, 3).Value = "Dpp4"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.4384066666666667
$ws.Cells.Item(7, 8).Value = 1.31522
$ws.Cells.Item(7, 9).Value = 0.9594354356895509
$ws.Cells.Item(7, 10).Value = 0.9594354356895509
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.711542
$ws.Cells.Item(7, 14).Value = 5.134626
$ws.Cells.Item(7, 15).Value = 0.2704366935850219
$ws.Cells.Item(7, 16).Value = 0.2704366935850219
$ws.Cells.Item(7, 17).Value = 0.7503514230799999
$ws.Cells.Item(7, 18).Value = 6.75316280772
$ws.Cells.Item(7, 19).Value = 0.2594665469361871
$ws.Cells.Item(7, 20).Value = 0.2594665469361871
